# Simulation ended early with too many parameter combinations; restart the
# run using fewer parameter combinations (wind_avg and wind_directions each
# collapsed to a single level/value), update the results-folder label and
# the "simulation started" timestamp, and move the active selection back
# to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# wind_avg: drop from 4 levels ("0, 0.2, 0.4, 0.6") down to a single value 0.2
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 0.2
$ws.Range("C27").HorizontalAlignment = -4131   # xlHAlignLeft

# wind_directions: drop from 2 levels ("all, up") down to a single value "all"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = "all"

# Update results folder name to reflect the restarted run (WIND -> light_WIND)
$ws.Range("A33").Value = "06-24-2014 - FP&SAV - light_WIND() - const NtoP - light_limitation - Scheffer vers"

# Update the simulation start time for the restarted run
$ws.Range("B1").Value = "6/24/2014 ~12:30PM"

# Restore the active selection to B2
$ws.Range("B2").Select()
